$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 23.70307803153992, 17.171875),
    @(37, 0.01, 3.178343534469604, 1.640625),
    @(50, 0.0016, 2.165600299835205, 1.40625),
    @(50, 0.0016, 12.8171751499176, 1.453125),
    @(50, 0.0016, 1.738581657409668, 0.53125),
    @(75, 0.0004, 74.63529062271118, 0.828125),
    @(39.44, 0.16, 30.566237449646, 20.5625),
    @(39.44, 0.16, 23.24444532394409, 18.25),
    @(75, 0.0004, 1.669373512268066, 0.53125),
    @(39.44, 0.16, 21.29177212715149, 15.109375),
    @(39.44, 0.16, 21.79749274253845, 14.71875),
    @(50, 0.0016, 2.104458570480347, 0.875),
    @(75, 0.0004, 1.485002517700195, 0.46875),
    @(33.33, 0.0036, 1.890873670578003, 0.9375),
    @(41.5, 0.04000000000000001, 7.52288556098938, 4.3125),
    @(41.5, 0.04000000000000001, 7.478296279907227, 3.484375),
    @(41.5, 0.04000000000000001, 9.29134726524353, 4.53125)
)

$startRow = 248
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
